# Apply updated benchmark metrics to the "Metrics" and "CLF_Report" worksheets,
# and refresh the Best_Params strings for Logistic_Regression, Random_Forest,
# Extra_Trees, Kernel_SVM and CatBoost on the "Metrics" sheet.

$wb = $excel.ActiveWorkbook

$metricsSheet = $wb.Worksheets.Item("Metrics")
$clfReportSheet = $wb.Worksheets.Item("CLF_Report")

# --- Updated Best_Params text (column M) on the Metrics sheet ---
$bestParamsValues = @{
    "M2" = '{''penalty'': ''l2'', ''solver'': ''liblinear'', ''C'': 353}'
    "M3" = '{''n_estimators'': 50, ''max_features'': ''auto'', ''criterion'': ''gini'', ''max_depth'': 4, ''min_samples_split'': 2}'
    "M4" = '{''n_estimators'': 600, ''max_features'': ''auto'', ''criterion'': ''gini'', ''max_depth'': 7, ''min_samples_split'': 5}'
    "M5" = '{''C'': 21, ''kernel'': ''rbf'', ''gamma'': 0.005686054187335557}'
    "M6" = '{''learning_rate'': 0.09232246106244313, ''depth'': 3, ''l2_leaf_reg'': 2, ''iterations'': 100}'
}

foreach ($cellRef in $bestParamsValues.Keys) {
    $metricsSheet.Range($cellRef).Value = $bestParamsValues[$cellRef]
}

# --- Updated numeric metrics (columns B-L, rows 2-9) on the Metrics sheet ---
$metricsValues = @{
    "B2" = 0.71
    "C2" = 0.64
    "D2" = 13
    "E2" = 8
    "F2" = 1
    "G2" = 11
    "H2" = 0.93
    "I2" = 0.54
    "J2" = 0.68
    "K2" = 0.89
    "L2" = 0.72
    "B3" = 0.83
    "C3" = 0.79
    "D3" = 18
    "E3" = 8
    "G3" = 6
    "H3" = 0.95
    "I3" = 0.75
    "J3" = 0.84
    "K3" = 0.89
    "B4" = 0.97
    "C4" = 0.76
    "D4" = 19
    "E4" = 6
    "F4" = 3
    "G4" = 5
    "H4" = 0.86
    "I4" = 0.79
    "J4" = 0.82
    "L4" = 0.73
    "C5" = 0.61
    "D5" = 12
    "E5" = 8
    "G5" = 12
    "H5" = 0.92
    "I5" = 0.5
    "J5" = 0.65
    "K5" = 0.89
    "L5" = 0.6899999999999999
    "B6" = 0.84
    "C6" = 0.76
    "D6" = 17
    "E6" = 8
    "G6" = 7
    "H6" = 0.9399999999999999
    "I6" = 0.71
    "J6" = 0.8100000000000001
    "K6" = 0.89
    "L6" = 0.8
    "C7" = 0.76
    "D7" = 18
    "E7" = 7
    "F7" = 2
    "G7" = 6
    "H7" = 0.9
    "I7" = 0.75
    "J7" = 0.82
    "K7" = 0.78
    "L7" = 0.76
    "B8" = 0.85
    "C8" = 0.79
    "D8" = 18
    "E8" = 8
    "F8" = 1
    "G8" = 6
    "H8" = 0.95
    "I8" = 0.75
    "J8" = 0.84
    "K8" = 0.89
    "L8" = 0.82
    "B9" = 1
    "C9" = 0.7
    "D9" = 20
    "E9" = 3
    "F9" = 6
    "G9" = 4
    "H9" = 0.77
    "I9" = 0.83
    "J9" = 0.8
    "L9" = 0.58
}

foreach ($cellRef in $metricsValues.Keys) {
    $metricsSheet.Range($cellRef).Value = $metricsValues[$cellRef]
}

# --- Updated numeric metrics (columns C-F, rows 2-41) on the CLF_Report sheet ---
$clfReportValues = @{
    "C2" = 0.42
    "D2" = 0.89
    "E2" = 0.57
    "F2" = 9
    "C3" = 0.93
    "D3" = 0.54
    "E3" = 0.68
    "F3" = 24
    "C4" = 0.64
    "D4" = 0.64
    "E4" = 0.64
    "F4" = 0.64
    "C5" = 0.67
    "D5" = 0.72
    "E5" = 0.63
    "F5" = 33
    "C6" = 0.79
    "D6" = 0.64
    "E6" = 0.65
    "F6" = 33
    "C7" = 0.57
    "D7" = 0.89
    "E7" = 0.7
    "F7" = 9
    "C8" = 0.95
    "D8" = 0.75
    "E8" = 0.84
    "F8" = 24
    "C9" = 0.79
    "D9" = 0.79
    "E9" = 0.79
    "F9" = 0.79
    "C10" = 0.76
    "E10" = 0.77
    "F10" = 33
    "D11" = 0.79
    "E11" = 0.8
    "F11" = 33
    "C12" = 0.57
    "D12" = 0.89
    "E12" = 0.7
    "F12" = 9
    "C13" = 0.95
    "D13" = 0.75
    "E13" = 0.84
    "F13" = 24
    "C14" = 0.79
    "D14" = 0.79
    "E14" = 0.79
    "F14" = 0.79
    "C15" = 0.76
    "E15" = 0.77
    "F15" = 33
    "D16" = 0.79
    "E16" = 0.8
    "F16" = 33
    "C17" = 0.4
    "D17" = 0.89
    "E17" = 0.55
    "F17" = 9
    "C18" = 0.92
    "D18" = 0.5
    "E18" = 0.65
    "F18" = 24
    "C19" = 0.61
    "D19" = 0.61
    "E19" = 0.61
    "F19" = 0.61
    "C20" = 0.66
    "D20" = 0.6899999999999999
    "E20" = 0.6
    "F20" = 33
    "C21" = 0.78
    "D21" = 0.61
    "E21" = 0.62
    "F21" = 33
    "C22" = 0.53
    "D22" = 0.89
    "E22" = 0.67
    "F22" = 9
    "C23" = 0.9399999999999999
    "D23" = 0.71
    "E23" = 0.8100000000000001
    "F23" = 24
    "C24" = 0.76
    "D24" = 0.76
    "E24" = 0.76
    "F24" = 0.76
    "C25" = 0.74
    "D25" = 0.8
    "E25" = 0.74
    "F25" = 33
    "C26" = 0.83
    "D26" = 0.76
    "E26" = 0.77
    "F26" = 33
    "C27" = 0.54
    "D27" = 0.78
    "E27" = 0.64
    "F27" = 9
    "C28" = 0.9
    "D28" = 0.75
    "E28" = 0.82
    "F28" = 24
    "C29" = 0.76
    "D29" = 0.76
    "E29" = 0.76
    "F29" = 0.76
    "C30" = 0.72
    "D30" = 0.76
    "E30" = 0.73
    "F30" = 33
    "C31" = 0.8
    "D31" = 0.76
    "E31" = 0.77
    "F31" = 33
    "D32" = 0.89
    "E32" = 0.7
    "F32" = 9
    "C33" = 0.95
    "D33" = 0.75
    "E33" = 0.84
    "F33" = 24
    "C34" = 0.79
    "D34" = 0.79
    "E34" = 0.79
    "F34" = 0.79
    "C35" = 0.76
    "D35" = 0.82
    "E35" = 0.77
    "F35" = 33
    "C36" = 0.84
    "D36" = 0.79
    "E36" = 0.8
    "F36" = 33
    "C37" = 0.43
    "E37" = 0.38
    "F37" = 9
    "C38" = 0.77
    "D38" = 0.83
    "E38" = 0.8
    "F38" = 24
    "C39" = 0.7
    "D39" = 0.7
    "E39" = 0.7
    "F39" = 0.7
    "C40" = 0.6
    "D40" = 0.58
    "E40" = 0.59
    "F40" = 33
    "C41" = 0.68
    "D41" = 0.7
    "E41" = 0.68
    "F41" = 33
}

foreach ($cellRef in $clfReportValues.Keys) {
    $clfReportSheet.Range($cellRef).Value = $clfReportValues[$cellRef]
}

Write-Host "Applied" $bestParamsValues.Count "Best_Params updates,"  $metricsValues.Count "Metrics updates, and" $clfReportValues.Count "CLF_Report updates."
